$wb = $excel.ActiveWorkbook

# ---- Sheet "2025" ----
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 0.076183
$ws.Range("E2").Value = 0.4079574890751533
$ws.Range("G2").Value = 0.2494892361374915
$ws.Range("I2").Value = 0.5195864764512724
$ws.Range("L2").Value = 0.5703349
$ws.Range("N2").Value = 11.85261458383828
$ws.Range("O2").Value = 2.574986507049832

# ---- Sheet "2030" ----
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 0.1296190864266778
$ws.Range("E2").Value = 0.3994426475317127
$ws.Range("I2").Value = 0.6524336873210169
$ws.Range("L2").Value = 0.328827
$ws.Range("N2").Value = 9.876658118432859
$ws.Range("O2").Value = 2.533576316016565

# ---- Sheet "2035" ----
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.04870720107345581
$ws.Range("B2").Value = 0.05335530511476796
$ws.Range("E2").Value = 0.2305989524682872
$ws.Range("I2").Value = 0.4607661397679059
$ws.Range("N2").Value = 9.047291122848666
$ws.Range("O2").Value = 6.165189263121612
